# Applies:
#  1) Bold the "Holly Dickson" name run
#  2) Bold the "Berufserfahrung" heading run
#  3) Bold the "Junior Animation Designer" heading run
#  4) Bold the "Animation – Praktikantin" heading run and retitle it to
#     "Praktikant im Bereich Animation"
#  5) Bold the "Bachelor of Fine Arts in Animation" heading run

$d = $word.ActiveDocument
$enDash = [char]8211

function Set-BoldOnText([string]$text) {
    $range = $d.Content
    $found = $range.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $range.Font.Bold = 1
    }
    return $found
}

Set-BoldOnText("Holly Dickson") | Out-Null
Set-BoldOnText("Berufserfahrung") | Out-Null
Set-BoldOnText("Junior Animation Designer") | Out-Null
Set-BoldOnText("Bachelor of Fine Arts in Animation") | Out-Null

# The "Animation – Praktikantin" run both gains bold and is retitled.
$oldTitle = "Animation " + $enDash + " Praktikantin"
$newTitle = "Praktikant im Bereich Animation"
$range = $d.Content
$found = $range.Find.Execute($oldTitle, $true, $false, $false, $false, $false, $true, 1, $false, $newTitle, 2)
if ($found) {
    $range.Font.Bold = 1
}
